$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# Force text format first so numeric-looking strings (e.g. "305.10") are not
# auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.654.72'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.276.97'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '305.10'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.68%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '96.42'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -2.62%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.499'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -3.13%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.46'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0789'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '18.35'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +3.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.118'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.72'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.651.26'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.282.80'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.777'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.561.24'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.02'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0896'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.99'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.11'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.90'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.13'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.54%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.45'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '4.01'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '25.09'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '166.11'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.61%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.04'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.25%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '32.97'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.75'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.97'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.55'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0689'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.74'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.69'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.68%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.996.26'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0279'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.77%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '18.17'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +3.72%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '9.94'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -3.24%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.06'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -6.23%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.76'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.95%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.85'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +4.33%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '53.56'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.506.42'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.94%  '
